# Applies the cryptocurrency price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-three (used inside the ShibaInu price, e.g. "0.0₃0880")
$sub3 = [string][char]0x2083

# Map of cell reference -> new text value.
$changes = @{}
$changes["D2"] = "43.933.19"
$changes["E2"] = "  -0.99%  "
$changes["D3"] = "2.196.62"
$changes["E3"] = "  -2.08%  "
$changes["E4"] = "  -0.19%  "
$changes["D5"] = "293.52"
$changes["E5"] = "  -4.16%  "
$changes["D6"] = "88.94"
$changes["E6"] = "  -4.77%  "
$changes["D7"] = "0.575"
$changes["E7"] = "  +0.56%  "
$changes["E8"] = "  -0.03%  "
$changes["D9"] = "0.480"
$changes["E9"] = "  -7.93%  "
$changes["D10"] = "32.43"
$changes["E10"] = "  -6.19%  "
$changes["D11"] = "0.0773"
$changes["E11"] = "  -4.42%  "
$changes["E12"] = "  -1.91%  "
$changes["D13"] = "6.72"
$changes["E13"] = "  -5.82%  "
$changes["B14"] = "WrappedliquidstakedEther2.0"
$changes["C14"] = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$changes["D14"] = "2.531.67"
$changes["E14"] = "  -2.28%  "
$changes["B15"] = "WrappedEther"
$changes["C15"] = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$changes["D15"] = "2.273.87"
$changes["E15"] = "  -3.66%  "
$changes["D16"] = "13.05"
$changes["E16"] = "  -3.67%  "
$changes["D17"] = "0.762"
$changes["E17"] = "  -8.53%  "
$changes["D18"] = "43.873.57"
$changes["E18"] = "  -0.39%  "
$changes["D19"] = "0.0" + $sub3 + "0880"
$changes["E19"] = "  -8.29%  "
$changes["D20"] = "5.82"
$changes["E20"] = "  -8.18%  "
$changes["D21"] = "10.74"
$changes["E21"] = "  -12.88%  "
$changes["D22"] = "62.98"
$changes["E22"] = "  -4.05%  "
$changes["D23"] = "228.85"
$changes["E23"] = "  -3.52%  "
$changes["D24"] = "2.75"
$changes["E24"] = "  -11.99%  "
$changes["E25"] = "  +0.02%  "
$changes["D26"] = "1.81"
$changes["E26"] = "  -8.25%  "
$changes["D27"] = "2.20"
$changes["E27"] = "  -0.34%  "
$changes["D28"] = "35.44"
$changes["E28"] = "  -8.54%  "
$changes["D29"] = "9.17"
$changes["E29"] = "  -6.30%  "
$changes["D30"] = "18.91"
$changes["E30"] = "  -5.38%  "
$changes["D31"] = "148.08"
$changes["E31"] = "  -3.36%  "
$changes["D32"] = "5.29"
$changes["E32"] = "  -10.91%  "
$changes["D33"] = "2.48"
$changes["E33"] = "  -5.93%  "
$changes["D34"] = "0.0731"
$changes["E34"] = "  -8.17%  "
$changes["D35"] = "0.116"
$changes["E35"] = "  -2.76%  "
$changes["D36"] = "2.85"
$changes["E36"] = "  -8.56%  "
$changes["E37"] = "  -7.27%  "
$changes["D38"] = "1.64"
$changes["E38"] = "  -7.50%  "
$changes["D39"] = "13.33"
$changes["E39"] = "  -9.15%  "
$changes["D40"] = "0.0279"
$changes["E40"] = "  -7.14%  "
$changes["D41"] = "3.06"
$changes["E41"] = "  -11.52%  "
$changes["B42"] = "FirstDigitalUSD"
$changes["C42"] = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$changes["D42"] = "1.01"
$changes["E42"] = "  -0.18%  "
$changes["B43"] = "RenderToken"
$changes["C43"] = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$changes["D43"] = "3.47"
$changes["E43"] = "  -8.86%  "
$changes["D44"] = "1.740.16"
$changes["E44"] = "  +0.01%  "
$changes["D45"] = "1.64"
$changes["E45"] = "  +0.28%  "
$changes["D46"] = "67.39"
$changes["E46"] = "  -1.45%  "
$changes["D47"] = "73.53"
$changes["E47"] = "  -9.02%  "
$changes["D48"] = "0.171"
$changes["E48"] = "  -10.58%  "
$changes["D49"] = "91.33"
$changes["D50"] = "2.415.07"
$changes["E50"] = "  -2.20%  "
$changes["D51"] = "7.42"
$changes["E51"] = "  -9.09%  "

# Write every value back as TEXT (matching the original inline-string cells) so that
# numeric-looking strings such as "293.52" or "2.20" are not reinterpreted as numbers
# and lose formatting (trailing zeros, grouping dots, the subscript digit, etc.).
foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = "Normal"
}
